# Fruta / hortaliza, semanal
# Applies the updated weekly values (Fecha, Volumen, Precio minimo/maximo/promedio,
# Origen, Precio $/Kg) to rows 2-31 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D=44575; J=61; K=8000; L=8000; M=8000; O='Provincia de Quillota'; P=500},
    @{Row=3; D=44407; J=45; K=5500; L=6000; M=5744; O='Provincia de Quillota'; P=359},
    @{Row=4; D=44341; J=51; K=5500; L=6000; M=5755; O='Provincia de Quillota'; P=360},
    @{Row=5; D=44582; J=52; K=7000; L=7000; M=7000; O='Provincia de Quillota'; P=438},
    @{Row=6; D=44715; J=70; K=5000; L=6000; M=5500; O='Provincia de Quillota'; P=344},
    @{Row=7; D=44330; J=120; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=8; D=44910; J=70; K=6000; L=7000; M=6500; O='Provincia de Quillota'; P=406},
    @{Row=9; D=44467; J=52; K=5000; L=6000; M=5500; O='Provincia de Quillota'; P=344},
    @{Row=10; D=44455; J=52; K=5000; L=6000; M=5500; O='Provincia de Quillota'; P=344},
    @{Row=11; D=44306; J=50; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=12; D=44691; J=61; K=6000; L=7000; M=6508; O='Provincia de Quillota'; P=407},
    @{Row=13; D=44371; J=34; K=5500; L=6000; M=5750; O='Provincia de Quillota'; P=359},
    @{Row=14; D=44358; J=52; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=15; D=44313; J=34; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=16; D=44698; J=34; K=6000; L=7000; M=6500; O='Provincia de Quillota'; P=406},
    @{Row=17; D=44308; J=70; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=18; D=44355; J=25; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=19; D=44403; J=43; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=20; D=44573; J=34; K=8000; L=8000; M=8000; O='Provincia de Quillota'; P=500},
    @{Row=21; D=44782; J=70; K=6000; L=6000; M=6000; O='Región Metropolitana'; P=375},
    @{Row=22; D=44328; J=160; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=23; D=44442; J=25; K=6000; L=7000; M=6480; O='Provincia de Quillota'; P=405},
    @{Row=24; D=44363; J=160; K=5500; L=6000; M=5750; O='Provincia de Quillota'; P=359},
    @{Row=25; D=44438; J=34; K=5000; L=6000; M=5500; O='Provincia de Quillota'; P=344},
    @{Row=26; D=44932; J=70; K=6000; L=7000; M=6500; O='Provincia de Quillota'; P=406},
    @{Row=27; D=44474; J=52; K=5000; L=6000; M=5500; O='Provincia de Quillota'; P=344},
    @{Row=28; D=44376; J=43; K=4500; L=5000; M=4756; O='Provincia de Quillota'; P=297},
    @{Row=29; D=44477; J=25; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=30; D=44350; J=25; K=6000; L=6000; M=6000; O='Provincia de Quillota'; P=375},
    @{Row=31; D=44589; J=52; K=8000; L=8000; M=8000; O='Provincia de Quillota'; P=500}
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D    # D: Fecha
    $ws.Cells.Item($r, 10).Value = $entry.J   # J: Volumen
    $ws.Cells.Item($r, 11).Value = $entry.K   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $entry.L   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $entry.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value = $entry.O   # O: Origen
    $ws.Cells.Item($r, 16).Value = $entry.P   # P: Precio $/Kg
}
